# Fixing geopoint in shared_table model
# Update header labels in the "model" sheet to append ".type" to the
# schema.properties.(latitude|longitude|altitude|accuracy) entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

$ws.Range("E1").Value = "schema.properties.latitude.type"
$ws.Range("F1").Value = "schema.properties.longitude.type"
$ws.Range("G1").Value = "schema.properties.altitude.type"
$ws.Range("H1").Value = "schema.properties.accuracy.type"
